# Update cryptos list values in column D (Price) and E (Volume 1h)
# D values are forced as text via a leading apostrophe (Excel text-prefix)
# because some look like numbers and would otherwise be auto-converted;
# the Style reset afterwards clears the quote-prefix flag so no extra
# cell style is left behind, matching the original unstyled text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'65.648.16"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -0.18%  "

$ws.Cells.Item(3, 4).Value = "'2.676.25"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -0.46%  "

$ws.Cells.Item(4, 5).Value = "  -0.02%  "

$ws.Cells.Item(5, 4).Value = "'600.72"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.19%  "

$ws.Cells.Item(6, 4).Value = "'156.93"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.26%  "

$ws.Cells.Item(7, 5).Value = "  -0.04%  "

$ws.Cells.Item(8, 4).Value = "'0.625"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +6.39%  "

$ws.Cells.Item(9, 4).Value = "'0.129"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.66%  "

$ws.Cells.Item(10, 4).Value = "'0.401"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.30%  "

$ws.Cells.Item(11, 4).Value = "'5.85"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.99%  "

$ws.Cells.Item(12, 5).Value = "  -0.28%  "

$ws.Cells.Item(13, 4).Value = "'29.37"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.74%  "

$ws.Cells.Item(14, 4).Value = "'0.0000198"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -0.38%  "

$ws.Cells.Item(15, 4).Value = "'3.157.53"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.72%  "

$ws.Cells.Item(16, 4).Value = "'65.519.43"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.17%  "

$ws.Cells.Item(17, 4).Value = "'2.671.72"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -0.82%  "

$ws.Cells.Item(18, 5).Value = "  +1.22%  "

$ws.Cells.Item(19, 5).Value = "  -1.72%  "

$ws.Cells.Item(20, 4).Value = "'7.58"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.01%  "

$ws.Cells.Item(21, 4).Value = "'351.79"
$ws.Cells.Item(21, 4).Style = "Normal"

$ws.Cells.Item(22, 5).Value = "  -0.05%  "

$ws.Cells.Item(23, 4).Value = "'69.61"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.80%  "

$ws.Cells.Item(24, 5).Value = "  +5.94%  "

$ws.Cells.Item(25, 4).Value = "'9.67"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.61%  "

$ws.Cells.Item(26, 4).Value = "'1.64"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +1.36%  "

$ws.Cells.Item(27, 5).Value = "  -0.71%  "

$ws.Cells.Item(28, 5).Value = "  -5.37%  "

$ws.Cells.Item(29, 4).Value = "'8.11"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.21%  "

$ws.Cells.Item(30, 5).Value = "  -0.17%  "

$ws.Cells.Item(31, 4).Value = "'533.46"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.95%  "

$ws.Cells.Item(32, 5).Value = "  -2.96%  "

$ws.Cells.Item(33, 5).Value = "  -2.02%  "

$ws.Cells.Item(34, 4).Value = "'6.46"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -2.11%  "

$ws.Cells.Item(35, 4).Value = "'5.50"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.07%  "

$ws.Cells.Item(36, 4).Value = "'0.424"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.71%  "

$ws.Cells.Item(37, 4).Value = "'20.50"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -1.15%  "

$ws.Cells.Item(38, 4).Value = "'0.999"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.02%  "

$ws.Cells.Item(39, 4).Value = "'158.20"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -3.07%  "

$ws.Cells.Item(40, 5).Value = "  -2.19%  "

$ws.Cells.Item(42, 4).Value = "'164.60"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -2.54%  "

$ws.Cells.Item(43, 5).Value = "  -0.28%  "

$ws.Cells.Item(44, 4).Value = "'2.33"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +2.79%  "

$ws.Cells.Item(45, 4).Value = "'0.0609"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.26%  "

$ws.Cells.Item(46, 4).Value = "'22.82"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.58%  "

$ws.Cells.Item(47, 5).Value = "  -1.87%  "

$ws.Cells.Item(48, 5).Value = "  -2.09%  "

$ws.Cells.Item(49, 4).Value = "'0.0₆0260"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +15.22%  "

$ws.Cells.Item(50, 5).Value = "  +2.66%  "

$ws.Cells.Item(51, 4).Value = "'20.14"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -4.79%  "
